$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 5, shifting existing rows (5..47) down to (6..48).
$ws.Rows.Item(5).Insert()

# Update the active selection to match the authored workbook state.
$ws.Range("E20").Select()
